# ---------------------------------------------------------------------------
# Applies the "Use Cases.docx" edits described in the commit diff.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-ParaIndex($doc, $text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $ptext = $p.Range.Text.TrimEnd()
        if ($ptext -eq $text) {
            return $i
        }
    }
    return -1
}

# Split the text inside a run range [start, start+oldSubLen) into its own
# run by replacing its text and then toggling Bold on/off (forces Word to
# break the surrounding run into three distinct <w:r> elements while
# leaving formatting unchanged). Returns the end offset of the new text.
function Split-Run($doc, $start, $oldSubLen, $newText) {
    $endOld = $start + $oldSubLen
    $r = $doc.Range($start, $endOld)
    $r.Text = $newText
    $newEnd = $start + $newText.Length
    $r2 = $doc.Range($start, $newEnd)
    $r2.Bold = 1
    $r2.Bold = 0
    return $newEnd
}

# ---------------------------------------------------------------------------
# 1) Remove "Program prompts log in or create new user" paragraph entirely.
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "Program prompts log in or create new user"
$d.Paragraphs.Item($idx).Range.Delete()

# ---------------------------------------------------------------------------
# 2) "User choses username & password" -> "User enters his/her username"
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "User choses username & password"
$d.Paragraphs.Item($idx).Range.Text = "User enters his/her username"

# ---------------------------------------------------------------------------
# 3) Remove "User opens application" and "Logs in" paragraphs from the
#    Delete-user section (the ones right before "Select Delete user option").
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "Select Delete user option"
$d.Paragraphs.Item($idx - 1).Range.Delete()
$idx = Find-ParaIndex $d "Select Delete user option"
$d.Paragraphs.Item($idx - 1).Range.Delete()

# ---------------------------------------------------------------------------
# 4) After "Select Delete user option" insert two new list paragraphs, then
#    append "d/denied" as a new run onto "Operation is complete".
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "Select Delete user option"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs.Item($idx).Range.Text = "User types in username to be deleted"

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs.Item($idx).Range.Text = "Program validates data"

$idx = Find-ParaIndex $d "Operation is complete"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter("d/denied")

# ---------------------------------------------------------------------------
# 5) "User is prompted to add server domain" -> three runs:
#    "User " + "selects 'Add Account'" + " to add server domain"
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "User is prompted to add server domain"
$p = $d.Paragraphs.Item($idx)
$pStart = $p.Range.Start
$s = $pStart + 5
$newWord = "selects " + [char]0x2018 + "Add Account" + [char]0x2019
Split-Run $d $s 11 $newWord

# ---------------------------------------------------------------------------
# 6) "User is prompted to select email address to be deleted" -> two runs:
#    "User " + "selects 'Remove Account'"
#    Then insert two new paragraphs after it.
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "User is prompted to select email address to be deleted"
$p = $d.Paragraphs.Item($idx)
$pStart = $p.Range.Start
$s = $pStart + 5
$oldLen = ("is prompted to select email address to be deleted").Length
$newWord = "selects " + [char]0x2018 + "Remove Account" + [char]0x2019
Split-Run $d $s $oldLen $newWord

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs.Item($idx).Range.Text = "User enters account to be deleted"

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx = $idx + 1
$d.Paragraphs.Item($idx).Range.Text = "Program validates data"

# ---------------------------------------------------------------------------
# 7) "System removes account with all data" + new run " or system denies operation"
# ---------------------------------------------------------------------------
$idx = Find-ParaIndex $d "System removes account with all data"
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertAfter(" or system denies operation")

# ---------------------------------------------------------------------------
# 8) "User is prompted is they're certain " ->
#    paragraph becomes "User selects " + "remove option" (bookmark retained)
#    new paragraph inserted after: "User is prompted if" + " they're certain "
# ---------------------------------------------------------------------------
$certainText = "User is prompted is they" + [char]0x2019 + "re certain"
$idx = Find-ParaIndex $d $certainText
$p = $d.Paragraphs.Item($idx)
$pStart = $p.Range.Start
$oldLen = ("is prompted is they" + [char]0x2019 + "re certain ").Length
$newWord = "selects "
Split-Run $d ($pStart + 5) $oldLen $newWord

$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$idx2 = $idx + 1
$newp = $d.Paragraphs.Item($idx2)
$newp.Range.Text = "User is prompted if" + " they" + [char]0x2019 + "re certain "

# now split "remove option" off from "User selects " run, and append "remove option"
$p = $d.Paragraphs.Item($idx)
$pStart = $p.Range.Start
$insertPoint = $pStart + 14
$r = $d.Range($insertPoint, $insertPoint)
$r.InsertAfter("remove option")
$r2start = $insertPoint
$r2end = $insertPoint + ("remove option").Length
$r2 = $d.Range($r2start, $r2end)
$r2.Bold = 1
$r2.Bold = 0

# split the new paragraph's text into two runs: "User is prompted if" + " they're certain "
$newp = $d.Paragraphs.Item($idx2)
$npStart = $newp.Range.Start
$splitAt = $npStart + ("User is prompted if").Length
$r3 = $d.Range($splitAt, $splitAt)
$r3.Bold = 1
$r3.Bold = 0
